$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: "Medians" label (new shared string)
$ws.Range("A18").Value = "Medians"

# Row 20: header row, same labels as row 3 (#year / whole / JanJul / AugSep)
$ws.Range("A20").Value = "      #   year      "
$ws.Range("B20").Value = "whole "
$ws.Range("C20").Value = "JanJul      "
$ws.Range("D20").Value = "AugSep"

# Row 21: 2013 medians
$ws.Range("A21").Value = 2013
$ws.Range("B21").Value = 6814.3919999999998
$ws.Range("C21").Value = 6147
$ws.Range("D21").Value = 5960.5
$ws.Range("F21").Formula = "=D21/C21"

# Row 22: 2014 medians
$ws.Range("A22").Value = 2014
$ws.Range("B22").Value = 7025.4040000000005
$ws.Range("C22").Value = 6301
$ws.Range("D22").Value = 6186
$ws.Range("F22").Formula = "=D22/C22"

# Row 23: 2015 medians
$ws.Range("A23").Value = 2015
$ws.Range("B23").Value = 7088.1270000000004
$ws.Range("C23").Value = 6530
$ws.Range("F23").Formula = "=(F22+F21)/2"
$ws.Range("D23").Formula = "=C23*F23"
$ws.Range("D23").Interior.Color = $ws.Range("D6").Interior.Color

# Row 25: average of the three median years
$ws.Range("C25").Formula = "=AVERAGE(C21:C23)"
$ws.Range("D25").Value = 6371.3501417178304
$ws.Range("F25").Formula = "=D25/C25"

# Row 31: carried-forward constants (mirrors row 14)
$ws.Range("C31").Value = 6889.3403333333335
$ws.Range("D31").Value = 6889.5730000000003
$ws.Range("F31").Formula = "=D31/C31"

# Selection left on F25, matching the saved workbook state
[void]$ws.Range("F25").Select()
